$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing bug: C57 end time was wrong (12:45 -> 00:45 next day) ---
$ws.Range("C57").Value2 = 1.03125

# --- Row 58: blank thick-bottom-border separator row before the new "Day 5" block ---
$ws.Range("A58:D58").RowHeight = 15.75

# --- Row 59: "Day 5" section header (copy formatting + labels from the "Day 4" header, row 47) ---
$ws.Range("A47:D47").Copy($ws.Range("A59:D59"))
$ws.Range("A59").Value = "Day 5"

# --- Row 60: first Day 5 entry (copy formatting from row 44, which has the same 2-line row height) ---
$ws.Range("B44:D44").Copy($ws.Range("B60:D60"))
$ws.Range("B60:D60").RowHeight = 30
$ws.Range("B60").Value = 0.39583333333333331
$ws.Range("C60").Value = 0.45833333333333331
$ws.Range("D60").Value = "Moved view settings to ClockViewStyle, implemented mark placement and ClockMarkView."

# --- Row 61: second Day 5 entry (copy formatting from row 45, default row height) ---
$ws.Range("B45:D45").Copy($ws.Range("B61:D61"))
$ws.Range("B61").Value = 0.46527777777777773
$ws.Range("C61").Value = 0.4826388888888889
$ws.Range("D61").Value = "Update test script for the new features."

# --- Row 62: third Day 5 entry (copy formatting from row 44, ht=30) ---
$ws.Range("B44:D44").Copy($ws.Range("B62:D62"))
$ws.Range("B62:D62").RowHeight = 30
$ws.Range("B62").Value = 0.5
$ws.Range("C62").Value = 0.5229166666666667
$ws.Range("D62").Value = "Debug mark placement. Marks almost placed correctly, right distance, but mirrored compared to the face."

# --- Row 63: fourth Day 5 entry (copy formatting from row 45, default row height); no "To" time yet ---
$ws.Range("B45:D45").Copy($ws.Range("B63:D63"))
$ws.Range("B63").Value = 0.52430555555555558
$ws.Range("C63").Clear()
$ws.Range("D63").Value = "I think the problem is in the face, not the marks."

# --- Keep the viewport in sync with the newly added rows ---
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("E68").Select()
